# Insert a new weekly price record as the new row 97, pushing the
# existing rows 97:151 down to 98:152 (dimension grows from A1:T151 to
# A1:T152).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 97 (shifts rows 97-151 down to 98-152).
$ws.Rows("97").Insert()

# Populate the newly inserted row 97 with the new data point.
$ws.Range("A97").Value = 5
$ws.Range("B97").Value = "Macroferia Regional de Talca"
$ws.Range("C97").Value = "Maule"
$ws.Range("D97").Value = 45264
$ws.Range("E97").Value = 7
$ws.Range("F97").Value = "Fruta"
$ws.Range("G97").Value = 100101
$ws.Range("H97").Value = "Berries"
$ws.Range("I97").Value = 100101001
$ws.Range("J97").Value = "Arándano (blue)"
$ws.Range("K97").Value = "Sin especificar"
$ws.Range("L97").Value = "Primera"
$ws.Range("M97").Value = 120
$ws.Range("N97").Value = 4600
$ws.Range("O97").Value = 4600
$ws.Range("P97").Value = 4600
$ws.Range("Q97").Value = "`$/bandeja 2 kilos"
$ws.Range("R97").Value = "Provincia de Curicó"
$ws.Range("S97").Value = 2300
$ws.Range("T97").Value = 2
